$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-16 changes from 45224 (2023-10-25)
# to 45233 (2023-11-03) serial date value.
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
